$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Address, $Val)
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Val
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '27.231.83'
Set-TextValue 'E2' '  +1.72%  '
Set-TextValue 'D3' '1.816.75'
Set-TextValue 'E3' '  +1.14%  '
Set-TextValue 'E4' '  +0.17%  '
Set-TextValue 'D5' '312.64'
Set-TextValue 'D6' '1.002'
Set-TextValue 'E6' '  +0.19%  '
Set-TextValue 'E7' '  +5.22%  '
Set-TextValue 'E8' '  +2.08%  '
Set-TextValue 'D9' '0.07397'
Set-TextValue 'E9' '  +0.48%  '
Set-TextValue 'E10' '  +1.67%  '
Set-TextValue 'E11' '  +0.04%  '
Set-TextValue 'D12' '1.818.12'
Set-TextValue 'E12' '  -3.83%  '
Set-TextValue 'D13' '6.655'
Set-TextValue 'E13' '  +0.48%  '
Set-TextValue 'D14' '5.394'
Set-TextValue 'E14' '  +2.68%  '
Set-TextValue 'D15' '0.07090'
Set-TextValue 'E15' '  +0.53%  '
Set-TextValue 'D16' '91.90'
Set-TextValue 'E16' '  -0.08%  '
Set-TextValue 'E17' '  +0.17%  '
Set-TextValue 'D18' '0.000008749'
Set-TextValue 'E18' '  +1.25%  '
Set-TextValue 'E19' '  +0.17%  '
Set-TextValue 'D20' '14.91'
Set-TextValue 'E20' '  +1.14%  '
Set-TextValue 'D21' '27.243.12'
Set-TextValue 'E21' '  +1.65%  '
Set-TextValue 'D22' '5.322'
Set-TextValue 'E22' '  +3.64%  '
Set-TextValue 'E23' '  +1.09%  '
Set-TextValue 'D24' '2.047.13'
Set-TextValue 'E24' '  -5.49%  '
Set-TextValue 'D25' '1.936'
Set-TextValue 'E25' '  -2.25%  '
Set-TextValue 'D26' '152.00'
Set-TextValue 'E26' '  +0.22%  '
Set-TextValue 'D27' '2.251'
Set-TextValue 'E27' '  +2.42%  '
Set-TextValue 'D28' '18.56'
Set-TextValue 'E28' '  +1.04%  '
Set-TextValue 'D29' '5.282'
Set-TextValue 'E29' '  +1.86%  '
Set-TextValue 'D30' '116.95'
Set-TextValue 'E30' '  +0.10%  '
Set-TextValue 'D31' '0.08884'
Set-TextValue 'E31' '  +1.23%  '
Set-TextValue 'D32' '0.7758'
Set-TextValue 'E32' '  +5.38%  '
Set-TextValue 'D34' '4.528'
Set-TextValue 'E34' '  +2.31%  '
Set-TextValue 'D35' '2.923'
Set-TextValue 'E35' '  +0.36%  '
Set-TextValue 'D36' '1.001'
Set-TextValue 'E36' '  +0.12%  '
Set-TextValue 'D37' '1.112'
Set-TextValue 'E37' '  +2.93%  '
Set-TextValue 'D38' '0.01962'
Set-TextValue 'D39' '0.05246'
Set-TextValue 'E39' '  +1.42%  '
Set-TextValue 'D40' '7.273'
Set-TextValue 'E40' '  +3.89%  '
Set-TextValue 'D41' '2.383'
Set-TextValue 'E41' '  +20.49%  '
Set-TextValue 'D42' '2.911'
Set-TextValue 'E42' '  +3.64%  '
Set-TextValue 'D43' '0.5292'
Set-TextValue 'E43' '  +1.52%  '
Set-TextValue 'D44' '0.1684'
Set-TextValue 'E44' '  +0.51%  '
Set-TextValue 'D45' '8.595'
Set-TextValue 'E45' '  +2.04%  '
Set-TextValue 'D46' '0.5048'
Set-TextValue 'E46' '  +0.96%  '
Set-TextValue 'D47' '10.47'
Set-TextValue 'E47' '  +1.31%  '
Set-TextValue 'D48' '105.12'
Set-TextValue 'E48' '  +0.08%  '
Set-TextValue 'B49' 'NEARProtocol'
Set-TextValue 'C49' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D49' '1.673'
Set-TextValue 'E49' '  +0.86%  '
Set-TextValue 'B50' 'PaxDollar'
Set-TextValue 'C50' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D50' '1.001'
Set-TextValue 'E50' '  +0.21%  '
Set-TextValue 'D51' '0.06326'
Set-TextValue 'E51' '  +0.27%  '
